$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ოზურგეთი")

# Update the E4:K4 values per the diff
$ws.Range("E4").Value = 14218
$ws.Range("F4").Value = 14418
$ws.Range("G4").Value = 14513
$ws.Range("H4").Value = 14643
$ws.Range("I4").Value = 14751
$ws.Range("J4").Value = 14991
$ws.Range("K4").Value = 15066

# Update the view: scroll so column B is the left-most visible column,
# and change the active selection to E4:K4
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("E4:K4").Select()
